$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 37748
$ws1.Range("F4").Value = 639
$ws1.Range("F5").Value = 782
$ws1.Range("F6").Value = 484
$ws1.Range("F7").Value = 371
$ws1.Range("F9").Value = 855
$ws1.Range("F11").Value = 726
$ws1.Range("F12").Value = 563
$ws1.Range("F13").Value = 61
$ws1.Range("F15").Value = 29
$ws1.Range("F16").Value = 662
$ws1.Range("F17").Value = 184
$ws1.Range("F18").Value = 473
$ws1.Range("F19").Value = 446
$ws1.Range("F20").Value = 1177
$ws1.Range("F22").Value = 846
$ws1.Range("F24").Value = 1036
$ws1.Range("F25").Value = 571
$ws1.Range("F27").Value = 1168
$ws1.Range("F29").Value = 797
$ws1.Range("F30").Value = 70
$ws1.Range("F31").Value = 1170

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 416
$ws2.Range("F10").Value = 13

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 643

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 643
$ws4.Range("F3").Value = 37748
$ws4.Range("F5").Value = 639
$ws4.Range("F6").Value = 782
$ws4.Range("F7").Value = 484
$ws4.Range("F9").Value = 371
$ws4.Range("F11").Value = 416
$ws4.Range("F15").Value = 855
$ws4.Range("F17").Value = 726
$ws4.Range("F18").Value = 563
$ws4.Range("F19").Value = 61
$ws4.Range("F24").Value = 13
$ws4.Range("F25").Value = 29
$ws4.Range("F27").Value = 662
$ws4.Range("F28").Value = 184
$ws4.Range("F29").Value = 473
$ws4.Range("F30").Value = 446
$ws4.Range("F31").Value = 1177
$ws4.Range("F33").Value = 846
$ws4.Range("F35").Value = 1036
$ws4.Range("F36").Value = 571
$ws4.Range("F38").Value = 1168
$ws4.Range("F41").Value = 797
$ws4.Range("F42").Value = 70
$ws4.Range("F43").Value = 1170
